{"js": "// Resume edit: the RPA/UiPath certification bullet changes from\n//   \"Successfully completed RPA Developer Foundation Diploma certification from UiPath\"\n// to\n//   \"Successfully completed RPA Advanced Diploma certification from UiPath\"\n//\n// i.e. the words \"Developer Foundation\" become \"Advanced\"; everything else\n// in that sentence (\" Diploma certification from \" / \"UiPath\") is unchanged.\n\nconst body = context.document.body;\n\n// --- Main content edit -----------------------------------------------\nconst target = body.search(\"Developer Foundation\", { matchCase: true });\ntarget.load(\"text\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  // This phrase only occurs once in the resume (the RPA/UiPath bullet).\n  target.items[0].insertText(\"Advanced\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Incidental run clean-up -------------------------------------------\n// The same save that introduced the edit above also re-serialized a couple\n// of other bullets, coalescing runs that carried identical formatting into\n// a single run (no visible text changes there). Reproduce that by\n// re-writing those ranges with their own (unchanged) text so adjacent\n// same-format runs collapse back into one.\n\n// \"course\" + \" completion certificate\" -> \"course completion certificate\"\n// (only the first \"Successful course completion certificate...\" bullet -\n// the Red Hat OpenShift one - actually has this run split).\nconst courseMatches = body.search(\"course completion certificate\", { matchCase: true });\ncourseMatches.load(\"text\");\nawait context.sync();\nif (courseMatches.items.length > 0) {\n  courseMatches.items[0].insertText(\"course completion certificate\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// \" Developer Certifi\" + \"cation Preparation\u201d\" -> \" Developer Certification Preparation\u201d\"\nconst certMatches = body.search(\"Developer Certification Preparation\\u201D\", { matchCase: true });\nawait context.sync();\nif (certMatches.items.length > 0) {\n  certMatches.items[0].insertText(\"Developer Certification Preparation\\u201D\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Certification line currently reads:\n#   \"Successfully completed RPA Developer Foundation Diploma certification from UiPath\"\n# Target:\n#   \"Successfully completed RPA Advanced Diploma certification from UiPath\"\n#\n# i.e. replace the words \"Developer Foundation\" with \"Advanced\" (the rest of\n# the sentence - \" Diploma certification from \" / \"UiPath\" - stays the same).\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"Developer Foundation\"\n$find.MatchCase = $true\n$find.Execute()\n\nif ($find.Found) {\n    $find.Parent.Text = \"Advanced\"\n}\n"}
